$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> (DAMSLTag, DialogAct)
$updates = @{
    8   = @("b",  "Acknowledge (Backchannel)")
    14  = @("sv", "Statement-opinion")
    23  = @("sv", "Statement-opinion")
    33  = @("sd", "Statement-non-opinion")
    58  = @("sv", "Statement-opinion")
    59  = @("aa", "Agree/Accept")
    64  = @("aa", "Agree/Accept")
    74  = @("sv", "Statement-opinion")
    92  = @("sv", "Statement-opinion")
    93  = @("sd", "Statement-non-opinion")
    97  = @("sv", "Statement-opinion")
    102 = @("sv", "Statement-opinion")
    105 = @("sd", "Statement-non-opinion")
    113 = @("sd", "Statement-non-opinion")
    130 = @("aa", "Agree/Accept")
    141 = @("sd", "Statement-non-opinion")
    144 = @("sd", "Statement-non-opinion")
    146 = @("sd", "Statement-non-opinion")
    148 = @("sv", "Statement-opinion")
    157 = @("sv", "Statement-opinion")
    160 = @("sd", "Statement-non-opinion")
    170 = @("aa", "Agree/Accept")
    171 = @("sv", "Statement-opinion")
    176 = @("sv", "Statement-opinion")
    185 = @("sv", "Statement-opinion")
    188 = @("sd", "Statement-non-opinion")
    193 = @("ba", "Appreciation")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
